$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: 2021 Take Your Child to Work Day ---
# Values are written in the same order the original authoring tool used, so that
# new entries land in xl/sharedStrings.xml in the same sequence.
$ws.Range("F2").Value2 = "Virtual Event"
$ws.Range("G2").Value2 = "https://bioinformatics.niaid.nih.gov/take-your-child-to-work"
$ws.Range("I2").Value2 = "No"
$ws.Range("J2").Value2 = "The NIH Office of Research Services hosts the annual Take Your Child to Work Day with an aim to inspire the next generation of NIH daughters and sons in grades 1-12 to explore career paths in science and public service at our nation’s medical research agency. 2021 is being held as a virtual event. Registration has closed, but resources provided by NIAID’s 3D Printing and Biovisualization Program are being made publicly available here."
$ws.Range("K2").Value2 = "3D,visualization,STEM"
$ws.Range("B2").Value2 = "2021 Take Your Child to Work Day"
$ws.Range("A2").Value2 = "NIH ALL STAFF LISTSERV"
$ws.Range("H2").Value2 = "Virtual Event"

$ws.Range("C2").Value2 = 44308
$ws.Range("C2").NumberFormat = "yyyy/m/dd"
$ws.Range("D2").Value2 = 44308
$ws.Range("D2").NumberFormat = "yyyy/m/dd"
$ws.Range("E2").Value2 = 1619096400

$ws.Rows.Item(2).RowHeight = 68

# --- Row 3: Advances in COVID-19 Prevention and Treatment workshop ---
$ws.Range("A3").Value2 = "STRUCTBIOLIG LISTSERV"
$ws.Range("B3").Value2 = "Advances in COVID-19 Prevention and Treatment Enabled by Structural Biology Research"
$ws.Range("F3").Value2 = "Virtual Workshop"
$ws.Range("G3").Value2 = "https://www.aps.anl.gov/sites/www.aps.anl.gov/files/APS-Uploads/WK9%20Agenda.pdf"
$ws.Range("H3").Value2 = "Virtual Event"
$ws.Range("I3").Value2 = "No"
$ws.Range("J3").Value2 = "Broadly, the workshop will present areas where structural biology research, including macromolecular crystallography and cryoelectron microscopy, intersects with in vivo, in vitro, and in silico studies of SARS-CoV-2 and COVID-19. More precisely, the topics will include (a) viral biology, (b) vaccine, therapeutic, and diagnostic antibody studies, and (c) small-molecule drug discovery as it relates to viral proteases and other viral proteins. In addition, as this year's events emphasize the need for a coordinated, long-term strategy to prevent future pandemics of zoonotic origin, a broader One Health perspective on viral pathogens will be presented."
$ws.Range("K3").Value2 = "structural biology,crystallography,SARS-CoV-2,drug discovery"

$ws.Range("C3").Value2 = 44327
$ws.Range("C3").NumberFormat = "yyyy/m/dd"
$ws.Range("D3").Value2 = 44328
$ws.Range("D3").NumberFormat = "yyyy/m/dd"
$ws.Range("E3").Value2 = 1620741600

$ws.Rows.Item(3).RowHeight = 102

# --- column H width (new column, best-fit like the C/D date columns) ---
$ws.Columns.Item(8).ColumnWidth = 10.83

# --- sheet view / selection state ---
$null = $ws.Range("J9").Select()
$excel.ActiveWindow.ScrollColumn = 7
$excel.ActiveWindow.ScrollRow = 1
